# Scheduled-runner update: refresh currentAveragePrice / LevePrice /
# LeveProfit columns (H:N) across the Leve-profit sheets with the latest
# market-board pull. Applied per-sheet, per-row.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param(
        $ws,
        [int]$Row,
        [hashtable]$Values,
        [array]$Clear
    )
    foreach ($col in $Values.Keys) {
        $ws.Cells.Item($Row, $col).Value = $Values[$col]
    }
    if ($Clear) {
        foreach ($col in $Clear) {
            $ws.Cells.Item($Row, $col).ClearContents()
        }
    }
}

# Column map: H=8 I=9 J=10 K=11 L=12 M=13 N=14

# ---------------------------------------------------------------- ALC ----
$ws = $wb.Worksheets.Item("ALC")

Set-Row $ws 32  @{8=14519.182; 9=2276.5; 10=17239.777; 11=2276.5; 12=17239.777; 13=-1950.5; 14=-17891.777} $null
Set-Row $ws 33  @{8=441.3125; 9=441.3125; 11=441.3125; 13=-212.3125} $null
Set-Row $ws 40  @{8=4448415; 9=3966.818; 10=37041036; 11=3966.818; 12=37041036; 13=-3791.818; 14=-37041386} $null
Set-Row $ws 64  @{8=4833; 9=4499; 11=4499; 13=-4251} $null
Set-Row $ws 67  @{8=4833; 9=4499; 11=4499; 13=-3641} $null
Set-Row $ws 98  @{8=1582.5454; 9=1582.5454; 10=0; 11=1582.5454; 12=0; 13=-84.54539999999997} @(14)
Set-Row $ws 112 @{8=82701.67999999999; 10=80266.234; 12=240798.702; 14=-243014.702} $null
Set-Row $ws 122 @{8=1582.5454; 9=1582.5454; 10=0; 11=4747.6362; 12=0; 13=-2297.6362} @(14)
Set-Row $ws 124 @{8=80000; 10=80000; 12=80000; 14=-89820} $null
Set-Row $ws 138 @{8=3507.4736; 9=2293.4119; 10=4490.2856; 11=6880.2357; 12=13470.8568; 13=-1740.2357; 14=-23750.8568} $null

# ---------------------------------------------------------------- ARM ----
$ws = $wb.Worksheets.Item("ARM")

Set-Row $ws 2   @{8=452787.66; 9=980823.5; 11=980823.5; 13=-980710.5} $null
Set-Row $ws 32  @{8=4244.213; 9=2718.8108; 11=2718.8108; 13=-2431.8108} $null
Set-Row $ws 45  @{8=1954.3334; 9=1954.3334; 11=1954.3334; 13=-1577.3334} $null
Set-Row $ws 116 @{8=452787.66; 9=980823.5; 11=980823.5; 13=-978529.5} $null
Set-Row $ws 122 @{8=1868.3914; 9=1483.8572; 11=4451.571599999999; 13=-2001.571599999999} $null
Set-Row $ws 132 @{8=2782704; 9=4352297.5; 11=13056892.5; 13=-13054362.5} $null
Set-Row $ws 140 @{8=62864.5; 10=62864.5; 12=62864.5; 14=-73224.5} $null

# ---------------------------------------------------------------- BSM ----
$ws = $wb.Worksheets.Item("BSM")

Set-Row $ws 3   @{8=452787.66; 9=980823.5; 11=980823.5; 13=-980709.5} $null
Set-Row $ws 82  @{8=9953; 9=9953; 11=9953; 13=-9570} $null
Set-Row $ws 85  @{8=9953; 9=9953; 11=9953; 13=-8627} $null
Set-Row $ws 134 @{8=11366066; 10=8750; 12=26250; 14=-31320} $null

# ---------------------------------------------------------------- CRP ----
$ws = $wb.Worksheets.Item("CRP")

Set-Row $ws 132 @{8=111113144; 9=111113144; 11=333339432; 13=-333336902} $null

# ---------------------------------------------------------------- CUL ----
$ws = $wb.Worksheets.Item("CUL")

Set-Row $ws 4   @{8=2669283.2; 10=2859914.2; 12=8579742.600000001; 14=-8579966.600000001} $null
Set-Row $ws 98  @{8=1622.375; 10=2300; 12=6900; 14=-9896} $null
Set-Row $ws 107 @{8=1966.1; 10=2427.75; 12=7283.25; 14=-11123.25} $null
Set-Row $ws 129 @{8=1197.619; 9=579.86664; 10=2742; 11=1739.59992; 12=8226; 13=3260.40008; 14=-18226} $null
Set-Row $ws 131 @{8=1331.6487; 9=804.3684; 10=1888.2222; 11=2413.1052; 12=5664.6666; 13=2626.8948; 14=-15744.6666} $null

# ---------------------------------------------------------------- GSM ----
$ws = $wb.Worksheets.Item("GSM")

Set-Row $ws 5   @{8=99999; 9=99999; 10=0; 11=99999; 12=0; 13=-99887} @(14)
Set-Row $ws 113 @{8=61524.176; 9=85337.086; 11=85337.086; 13=-83167.086} $null
Set-Row $ws 122 @{8=179452.72; 9=302792.25; 10=15000; 11=908376.75; 12=45000; 13=-905926.75; 14=-49900} $null
Set-Row $ws 126 @{8=3878.175; 9=3861.4856; 11=11584.4568; 13=-9114.4568} $null
Set-Row $ws 134 @{8=19666.666; 10=19666.666; 12=58999.99800000001; 14=-64069.99800000001} $null

# ---------------------------------------------------------------- LTW ----
$ws = $wb.Worksheets.Item("LTW")

Set-Row $ws 46  @{8=966.25; 9=818.4286; 11=818.4286; 13=-630.4286} $null
Set-Row $ws 61  @{8=2454.9375; 9=2213.7693; 11=2213.7693; 13=-2011.7693} $null
Set-Row $ws 93  @{8=500; 9=500; 10=0; 11=500; 12=0; 13=748} @(14)
Set-Row $ws 111 @{8=69999; 10=69999; 12=69999; 14=-78179} $null
Set-Row $ws 113 @{8=2454.9375; 9=2213.7693; 11=2213.7693; 13=-43.76929999999993} $null
Set-Row $ws 122 @{8=12221; 9=12221; 11=36663; 13=-34213} $null
Set-Row $ws 132 @{8=34301900; 9=34301900; 11=102905700; 13=-102903170} $null

# ---------------------------------------------------------------- WVR ----
$ws = $wb.Worksheets.Item("WVR")

Set-Row $ws 16  @{8=75000; 10=75000; 12=75000; 14=-75584} $null
Set-Row $ws 96  @{8=2539.2; 9=2299.8; 10=2778.6; 11=2299.8; 12=2778.6; 13=-926.8000000000002; 14=-5524.6} $null
Set-Row $ws 132 @{8=29423654; 10=27474; 12=82422; 14=-87482} $null
